$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value2 = $cell.Value2 + 1
}
